$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "Metering" column header to "Port-Range"
$ws.Range("C1").Value = "Port-Range"

# Update the selected cell to reflect the edited cell
$ws.Range("C2").Select()
